$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 26; this shifts existing rows 26..67 down to 27..68,
# preserving their values/formatting (matches the diff's row-shift pattern).
$ws.Rows(26).Insert()

# Populate the newly inserted row 26 with the new record's data.
$ws.Range("A26").Value = 4
$ws.Range("B26").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C26").Value = "Los Lagos"
$ws.Range("D26").Value = 44757
$ws.Range("E26").Value = 10
$ws.Range("F26").Value = 100112043
$ws.Range("G26").Value = "Pepino dulce"
$ws.Range("H26").Value = "Cultivar IV Región"
$ws.Range("I26").Value = "Primera"
$ws.Range("J26").Value = 50
$ws.Range("K26").Value = 18000
$ws.Range("L26").Value = 18000
$ws.Range("M26").Value = 18000
$ws.Range("N26").Value = "$/bandeja 18 kilos"
$ws.Range("O26").Value = "Provincia de Limarí"
$ws.Range("P26").Value = 1000
$ws.Range("Q26").Value = 18
$ws.Range("R26").Value = "Hortaliza"
